$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.169133186340332
$ws.Range("B1").Value = 2.377488851547241
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.382237434387207
$ws.Range("E1").Value = 1.211309909820557
